$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text so numeric-looking strings
# (trailing zeros, multi-dot thousand separators, etc.) are preserved exactly.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.606.10'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.923.08'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.46'
$ws.Range('E5').Value = '  +3.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2907'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06878'
$ws.Range('E9').Value = '  +4.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '105.40'
$ws.Range('E10').Value = '  -1.45%  '
$ws.Range('E11').Value = '  -3.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.919.48'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07719'
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.350'
$ws.Range('E14').Value = '  +3.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6703'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '288.28'
$ws.Range('E16').Value = '  -6.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.614.33'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007643'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.549'
$ws.Range('E21').Value = '  +4.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.172.59'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.444'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.505'
$ws.Range('E25').Value = '  +2.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.84'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.75'
$ws.Range('E27').Value = '  +2.32%  '
$ws.Range('E28').Value = '  +4.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1071'
$ws.Range('E29').Value = '  -3.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.409'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.183'
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.056'
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05022'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7339'
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02074'
$ws.Range('E36').Value = '  +5.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9996'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.731'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.687'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.044'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '111.26'
$ws.Range('E41').Value = '  +3.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4435'
$ws.Range('E42').Value = '  +6.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8750'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.895'
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.76'
$ws.Range('E46').Value = '  -4.47%  '
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.393'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1246'
$ws.Range('E49').Value = '  +3.22%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '47.80'
$ws.Range('E50').Value = '  +12.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.04'
$ws.Range('E51').Value = '  +0.44%  '
